$d = $word.ActiveDocument

# --- Step 1: append 5 new paragraphs after the final "</html>" paragraph ---
# (this leaves the "_GoBack" bookmark sitting on the old last paragraph,
#  right after the "</html>" text)
for ($i = 0; $i -lt 5; $i++) {
    $r = $d.Content
    $r.Collapse(0)
    $r.InsertParagraphAfter()
}

# --- Step 2: put the new URL text (plus one throwaway placeholder char) into
#             the new, still-empty, final paragraph ---
$lastPara = $d.Paragraphs.Last
$insertRange = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)
$insertRange.InsertAfter("https://github.com/Bigchance64/peetro-broker-farsi.gitX")

# --- Step 3: move the "_GoBack" bookmark onto the new final paragraph, right
#             after the real text (i.e. right before the placeholder char) ---
$lastPara = $d.Paragraphs.Last
$pos = $lastPara.Range.End - 2
$bmRange = $d.Range($pos, $pos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# --- Step 4: remove the placeholder char now that the bookmark is anchored
#             in place (this also removes the bookmark from the old
#             "</html>" paragraph, since "_GoBack" only ever has one
#             location) ---
$lastPara = $d.Paragraphs.Last
$placeholderRange = $d.Range($lastPara.Range.End - 2, $lastPara.Range.End - 1)
$placeholderRange.Delete()
